$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet currently has rows 2-15 holding line1..line6 then extr1..extr8.
# Two new rows (line7, line8) need to be inserted right after line6 (row 7),
# i.e. before the current row 8 (extr1). That pushes extr1..extr8 from rows
# 8-15 down to rows 10-17.
#
# Shift the extr1..extr8 block down by two rows, working bottom-to-top so
# that source rows are not overwritten before they are copied. Also copy the
# id-column (A) formatting along with the values so the new rows keep the
# existing look (bold font + border) instead of picking up a blank style.
for ($r = 15; $r -ge 8; $r--) {
    $destRow = $r + 2
    $ws.Cells.Item($destRow, 1).Value = $ws.Cells.Item($r, 1).Value2
    $ws.Cells.Item($destRow, 2).Value = $ws.Cells.Item($r, 2).Value2
    $ws.Cells.Item($destRow, 3).Value = $ws.Cells.Item($r, 3).Value2
    $ws.Cells.Item($destRow, 4).Value = $ws.Cells.Item($r, 4).Value2
    $ws.Cells.Item($destRow, 5).Value = [bool]($ws.Cells.Item($r, 5).Value2)

    $ws.Range("A$r").Copy() | Out-Null
    $ws.Range("A$destRow").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
}
$excel.CutCopyMode = 0

# Fill in the two new rows: line7 (row 8) and line8 (row 9).
$ws.Cells.Item(8, 1).Value = 6
$ws.Cells.Item(8, 2).Value = "line7"
$ws.Cells.Item(8, 3).Value = 14
$ws.Cells.Item(8, 4).Value = 11
$ws.Cells.Item(8, 5).Value = $true

$ws.Cells.Item(9, 1).Value = 7
$ws.Cells.Item(9, 2).Value = "line8"
$ws.Cells.Item(9, 3).Value = 16
$ws.Cells.Item(9, 4).Value = 9
$ws.Cells.Item(9, 5).Value = $true

# Renumber the id column for the rows that now hold extr1..extr8 (rows 10-17).
for ($r = 10; $r -le 17; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 2
}

# Apply the remaining data changes that landed on the shifted extr rows.
# Row 10 = extr1: in_service False -> True
$ws.Cells.Item(10, 5).Value = $true
# Row 11 = extr2: unchanged
# Row 12 = extr3: from_bus 9 -> 10
$ws.Cells.Item(12, 3).Value = 10
# Row 13 = extr4: in_service True -> False
$ws.Cells.Item(13, 5).Value = $false
# Row 14 = extr5: in_service True -> False
$ws.Cells.Item(14, 5).Value = $false
# Row 15 = extr6: unchanged
# Row 16 = extr7: unchanged
# Row 17 = extr8: in_service False -> True
$ws.Cells.Item(17, 5).Value = $true
